# Update LR-pair edge statistics per revised NATMI calculation (Dr Hou advice)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 13.739149
$ws.Range("H2").Value = 41.217447
$ws.Range("I2").Value = 0.6130043224686931
$ws.Range("J2").Value = 0.6130043224686931
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 297.8183156666666
$ws.Range("N2").Value = 893.4549469999999
$ws.Range("O2").Value = 0.8852156413092672
$ws.Range("P2").Value = 0.8852156413092673
$ws.Range("Q2").Value = 4091.770213873367
$ws.Range("R2").Value = 36825.93192486031
$ws.Range("S2").Value = 0.542641014439477
$ws.Range("T2").Value = 0.5426410144394771
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 13.739149
$ws.Range("H3").Value = 41.217447
$ws.Range("I3").Value = 0.6130043224686931
$ws.Range("J3").Value = 0.6130043224686931
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 24.34034433333333
$ws.Range("N3").Value = 73.021033
$ws.Range("O3").Value = 0.07234764413494278
$ws.Range("P3").Value = 0.0723476441349428
$ws.Range("Q3").Value = 334.4156175069724
$ws.Range("R3").Value = 3009.740557562751
$ws.Range("S3").Value = 0.04434941857514672
$ws.Range("T3").Value = 0.04434941857514672
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 13.739149
$ws.Range("H4").Value = 41.217447
$ws.Range("I4").Value = 0.6130043224686931
$ws.Range("J4").Value = 0.6130043224686931
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.277234
$ws.Range("N4").Value = 42.831702
$ws.Range("O4").Value = 0.04243671455578994
$ws.Range("P4").Value = 0.04243671455578994
$ws.Range("Q4").Value = 196.157045233866
$ws.Range("R4").Value = 1765.413407104794
$ws.Range("S4").Value = 0.02601388945406934
$ws.Range("T4").Value = 0.02601388945406934
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 4.034036666666666
$ws.Range("H5").Value = 12.10211
$ws.Range("I5").Value = 0.1799879973398545
$ws.Range("J5").Value = 0.1799879973398545
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 297.8183156666666
$ws.Range("N5").Value = 893.4549469999999
$ws.Range("O5").Value = 0.8852156413092672
$ws.Range("P5").Value = 0.8852156413092673
$ws.Range("Q5").Value = 1201.410005404241
$ws.Range("R5").Value = 10812.69004863817
$ws.Range("S5").Value = 0.15932819049317
$ws.Range("T5").Value = 0.15932819049317
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4.034036666666666
$ws.Range("H6").Value = 12.10211
$ws.Range("I6").Value = 0.1799879973398545
$ws.Range("J6").Value = 0.1799879973398545
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 24.34034433333333
$ws.Range("N6").Value = 73.021033
$ws.Range("O6").Value = 0.07234764413494278
$ws.Range("P6").Value = 0.0723476441349428
$ws.Range("Q6").Value = 98.18984151995889
$ws.Range("R6").Value = 883.70857367963
$ws.Range("S6").Value = 0.01302170758010482
$ws.Range("T6").Value = 0.01302170758010483
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4.034036666666666
$ws.Range("H7").Value = 12.10211
$ws.Range("I7").Value = 0.1799879973398545
$ws.Range("J7").Value = 0.1799879973398545
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.277234
$ws.Range("N7").Value = 42.831702
$ws.Range("O7").Value = 0.04243671455578994
$ws.Range("P7").Value = 0.04243671455578994
$ws.Range("Q7").Value = 57.59488545457999
$ws.Range("R7").Value = 518.3539690912199
$ws.Range("S7").Value = 0.007638099266579684
$ws.Range("T7").Value = 0.007638099266579685
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.639623666666666
$ws.Range("H8").Value = 13.918871
$ws.Range("I8").Value = 0.2070076801914524
$ws.Range("J8").Value = 0.2070076801914524
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 297.8183156666666
$ws.Range("N8").Value = 893.4549469999999
$ws.Range("O8").Value = 0.8852156413092672
$ws.Range("P8").Value = 0.8852156413092673
$ws.Range("Q8").Value = 1381.764905733871
$ws.Range("R8").Value = 12435.88415160484
$ws.Range("S8").Value = 0.1832464363766202
$ws.Range("T8").Value = 0.1832464363766202
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.639623666666666
$ws.Range("H9").Value = 13.918871
$ws.Range("I9").Value = 0.2070076801914524
$ws.Range("J9").Value = 0.2070076801914524
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 24.34034433333333
$ws.Range("N9").Value = 73.021033
$ws.Range("O9").Value = 0.07234764413494278
$ws.Range("P9").Value = 0.0723476441349428
$ws.Range("Q9").Value = 112.9300376237492
$ws.Range("R9").Value = 1016.370338613743
$ws.Range("S9").Value = 0.01497651797969124
$ws.Range("T9").Value = 0.01497651797969125
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.639623666666666
$ws.Range("H10").Value = 13.918871
$ws.Range("I10").Value = 0.2070076801914524
$ws.Range("J10").Value = 0.2070076801914524
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.277234
$ws.Range("N10").Value = 42.831702
$ws.Range("O10").Value = 0.04243671455578994
$ws.Range("P10").Value = 0.04243671455578994
$ws.Range("Q10").Value = 66.240992760938
$ws.Range("R10").Value = 596.1689348484419
$ws.Range("S10").Value = 0.008784725835140917
$ws.Range("T10").Value = 0.008784725835140918
